$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

# Update Runmode column: rows 2-7 flip from Yes to No, rows 8-12 flip from No to Yes
$ws.Range("E2:E7").Value = "No"
$ws.Range("E8:E12").Value = "Yes"

# Update the active selection to match the new state
$ws.Range("E8:E12").Select()
